$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2272727272727273
$ws.Range("C2").Value = 0.5454545454545454
$ws.Range("J2").Value = 0.04545454545454546
$ws.Range("P2").Value = 0.1818181818181818
$ws.Range("J3").Value = 0.08333333333333333
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("P4").Value = 0.6
$ws.Range("S4").Value = 0.4
$ws.Range("D6").Value = 0.05555555555555555
$ws.Range("F6").Value = 0.05555555555555555
$ws.Range("J6").Value = 0.3888888888888889
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.2777777777777778
$ws.Range("B7").Value = 0.2222222222222222
$ws.Range("D7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.1111111111111111
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.4444444444444444
$ws.Range("B8").Value = 0.1428571428571428
$ws.Range("D8").Value = 0.03571428571428571
$ws.Range("F8").Value = 0.07142857142857142
$ws.Range("J8").Value = 0.1071428571428571
$ws.Range("O8").Value = 0.03571428571428571
$ws.Range("Q8").Value = 0.1428571428571428
$ws.Range("R8").Value = 0.1428571428571428
$ws.Range("S8").Value = 0.3214285714285715
$ws.Range("J9").Value = 0.1818181818181818
$ws.Range("Q9").Value = 0.09090909090909091
$ws.Range("S9").Value = 0.7272727272727273
$ws.Range("B10").Value = 0.0873015873015873
$ws.Range("D10").Value = 0.01587301587301587
$ws.Range("F10").Value = 0.04761904761904762
$ws.Range("J10").Value = 0.1031746031746032
$ws.Range("O10").Value = 0.01587301587301587
$ws.Range("Q10").Value = 0.246031746031746
$ws.Range("R10").Value = 0.07936507936507936
$ws.Range("S10").Value = 0.4047619047619048
$ws.Range("G11").Value = 0.1052631578947368
$ws.Range("J11").Value = 0.1052631578947368
$ws.Range("K11").Value = 0.2105263157894737
$ws.Range("L11").Value = 0.5263157894736842
$ws.Range("S11").Value = 0.05263157894736842
$ws.Range("G12").Value = 0.3
$ws.Range("J12").Value = 0.5
$ws.Range("S12").Value = 0.2
$ws.Range("G13").Value = 1
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.4285714285714285
$ws.Range("O15").Value = 0.07142857142857142
$ws.Range("S15").Value = 0.2142857142857143
$ws.Range("H16").Value = 0.1875
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.4375
$ws.Range("K16").Value = 0.125
$ws.Range("O16").Value = 0.0625
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.025
$ws.Range("H17").Value = 0.175
$ws.Range("I17").Value = 0.075
$ws.Range("J17").Value = 0.475
$ws.Range("K17").Value = 0.05
$ws.Range("O17").Value = 0.05
$ws.Range("S17").Value = 0.15
$ws.Range("F18").Value = 0.06666666666666667
$ws.Range("H18").Value = 0.1333333333333333
$ws.Range("J18").Value = 0.6666666666666666
$ws.Range("S18").Value = 0.1333333333333333
$ws.Range("F19").Value = 0.0423728813559322
$ws.Range("H19").Value = 0.1186440677966102
$ws.Range("I19").Value = 0.0423728813559322
$ws.Range("J19").Value = 0.4152542372881356
$ws.Range("K19").Value = 0.09322033898305085
$ws.Range("M19").Value = 0.03389830508474576
$ws.Range("O19").Value = 0.05084745762711865
$ws.Range("S19").Value = 0.2033898305084746
